$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values scraped on Sun Nov 10 21:09:36 UTC 2024

# Row 2
$ws.Range("D2").Value = "78.856.66"
$ws.Range("E2").Value = "  +3.41%  "

# Row 3
$ws.Range("D3").Value = "3.113.34"
$ws.Range("E3").Value = "  +1.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.255"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +21.78%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.24%  "

# Row 10
$ws.Range("D10").Value = "3.110.54"
$ws.Range("E10").Value = "  +0.94%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.561"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +25.83%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.163"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.48%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +21.71%  "

# Row 14
$ws.Range("D14").Value = "3.693.87"
$ws.Range("E14").Value = "  +1.27%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.94%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.70%  "

# Row 17
$ws.Range("D17").Value = "78.985.98"
$ws.Range("E17").Value = "  +3.68%  "

# Row 18
$ws.Range("D18").Value = "3.130.92"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.52%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "426.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.55%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.66%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.07%  "

# Row 25
$ws.Range("D25").Value = "3.297.96"
$ws.Range("E25").Value = "  +1.85%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "74.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.67%  "

# Row 29
$ws.Range("E29").Value = "  -0.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000116"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.46%  "

# Row 31
$ws.Range("E31").Value = "  +0.87%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "532.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.87%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.53%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.144"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.37%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.119"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.86%  "

# Row 39
$ws.Range("E39").Value = "  -0.11%  "

# Row 40
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "20.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.04%  "

# Row 41
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.392"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.42%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "

# Row 43
$ws.Range("E43").Value = "  -0.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "183.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.11%  "

# Row 48
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.80%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.761"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.92%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.39%  "
